$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 154.54546
$ws.Range("I4").Value = 154.54546
$ws.Range("K4").Value = 154.54546
$ws.Range("M4").Value = -40.54545999999999
$ws.Range("H19").Value = 1716.7142
$ws.Range("J19").Value = 2244.8
$ws.Range("L19").Value = 2244.8
$ws.Range("N19").Value = -2594.8
$ws.Range("H28").Value = 44901.13
$ws.Range("I28").Value = 67098.13
$ws.Range("J28").Value = 3281.75
$ws.Range("K28").Value = 67098.13
$ws.Range("L28").Value = 3281.75
$ws.Range("M28").Value = -66613.13
$ws.Range("N28").Value = -4251.75
$ws.Range("H40").Value = 5731.8
$ws.Range("J40").Value = 7283.4287
$ws.Range("L40").Value = 7283.4287
$ws.Range("N40").Value = -7633.4287
$ws.Range("H101").Value = 485.2
$ws.Range("I101").Value = 485.2
$ws.Range("K101").Value = 1455.6
$ws.Range("M101").Value = 166.4000000000001
$ws.Range("H112").Value = 3251.7673
$ws.Range("J112").Value = 3311.6904
$ws.Range("L112").Value = 9935.0712
$ws.Range("N112").Value = -12151.0712
$ws.Range("H113").Value = 4751.25
$ws.Range("I113").Value = 5502.5
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 5502.5
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -2248.5
$ws.Range("N113").Value = -10508
$ws.Range("H132").Value = 2129643.2
$ws.Range("I132").Value = 2632664
$ws.Range("K132").Value = 7897992
$ws.Range("M132").Value = -7895462
$ws.Range("H137").Value = 6754.524
$ws.Range("I137").Value = 4513.857
$ws.Range("J137").Value = 7874.857
$ws.Range("K137").Value = 13541.571
$ws.Range("L137").Value = 23624.571
$ws.Range("M137").Value = -10991.571
$ws.Range("N137").Value = -28724.571
$ws.Range("H138").Value = 3802.2158
$ws.Range("I138").Value = 1615.5116
$ws.Range("J138").Value = 5891.7334
$ws.Range("K138").Value = 4846.5348
$ws.Range("L138").Value = 17675.2002
$ws.Range("M138").Value = 293.4651999999996
$ws.Range("N138").Value = -27955.2002
$ws.Range("H141").Value = 1696.4736
$ws.Range("I141").Value = 1696.4736
$ws.Range("K141").Value = 5089.4208
$ws.Range("M141").Value = 90.57920000000013

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3453.422
$ws.Range("I32").Value = 2721.3948
$ws.Range("K32").Value = 2721.3948
$ws.Range("M32").Value = -2434.3948
$ws.Range("H45").Value = 8299.200000000001
$ws.Range("I45").Value = 5998.6665
$ws.Range("K45").Value = 5998.6665
$ws.Range("M45").Value = -5621.6665
$ws.Range("H61").Value = 2127.608
$ws.Range("I61").Value = 1459.2727
$ws.Range("K61").Value = 1459.2727
$ws.Range("M61").Value = -1247.2727
$ws.Range("H74").Value = 1410.8889
$ws.Range("I74").Value = 1381.6
$ws.Range("J74").Value = 1645.2
$ws.Range("K74").Value = 1381.6
$ws.Range("L74").Value = 1645.2
$ws.Range("M74").Value = -507.5999999999999
$ws.Range("N74").Value = -3393.2
$ws.Range("H77").Value = 1410.8889
$ws.Range("I77").Value = 1381.6
$ws.Range("J77").Value = 1645.2
$ws.Range("K77").Value = 6908
$ws.Range("L77").Value = 8226
$ws.Range("M77").Value = -2540
$ws.Range("N77").Value = -16962
$ws.Range("H97").Value = 1033.6
$ws.Range("I97").Value = 1159.2778
$ws.Range("J97").Value = 530.8889
$ws.Range("K97").Value = 1159.2778
$ws.Range("L97").Value = 530.8889
$ws.Range("M97").Value = -663.2778000000001
$ws.Range("N97").Value = -1522.8889
$ws.Range("H132").Value = 2635.9507
$ws.Range("I132").Value = 875.17645
$ws.Range("K132").Value = 2625.52935
$ws.Range("M132").Value = -95.52935000000025
$ws.Range("H136").Value = 2127.608
$ws.Range("I136").Value = 1459.2727
$ws.Range("K136").Value = 4377.8181
$ws.Range("M136").Value = -1827.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1022.75
$ws.Range("I94").Value = 424.4
$ws.Range("K94").Value = 424.4
$ws.Range("M94").Value = 26.60000000000002
$ws.Range("H107").Value = 1880.3846
$ws.Range("I107").Value = 1821.3043
$ws.Range("K107").Value = 1821.3043
$ws.Range("M107").Value = 98.69569999999999
$ws.Range("H134").Value = 26423.977
$ws.Range("I134").Value = 2632.1516
$ws.Range("K134").Value = 7896.4548
$ws.Range("M134").Value = -5361.4548

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1674339.4
$ws.Range("J31").Value = 13000
$ws.Range("L31").Value = 13000
$ws.Range("N31").Value = -13590
$ws.Range("H34").Value = 1674339.4
$ws.Range("J34").Value = 13000
$ws.Range("L34").Value = 13000
$ws.Range("N34").Value = -13404
$ws.Range("H105").Value = 2397.2856
$ws.Range("I105").Value = 2296.8333
$ws.Range("K105").Value = 2296.8333
$ws.Range("M105").Value = -549.8332999999998
$ws.Range("H107").Value = 378.4
$ws.Range("I107").Value = 652.25
$ws.Range("J107").Value = 309.9375
$ws.Range("K107").Value = 652.25
$ws.Range("L107").Value = 309.9375
$ws.Range("M107").Value = 1267.75
$ws.Range("N107").Value = -4149.9375
$ws.Range("H132").Value = 3166.3333
$ws.Range("I132").Value = 2268.2327
$ws.Range("J132").Value = 5924.7856
$ws.Range("K132").Value = 6804.6981
$ws.Range("L132").Value = 17774.3568
$ws.Range("M132").Value = -4274.6981
$ws.Range("N132").Value = -22834.3568
$ws.Range("H134").Value = 253686.3
$ws.Range("I134").Value = 2511.36
$ws.Range("K134").Value = 7534.08
$ws.Range("M134").Value = -4999.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 212
$ws.Range("I23").Value = 267.5
$ws.Range("J23").Value = 175
$ws.Range("K23").Value = 802.5
$ws.Range("L23").Value = 525
$ws.Range("M23").Value = -567.5
$ws.Range("N23").Value = -995
$ws.Range("H34").Value = 56465.95
$ws.Range("I34").Value = 146.55556
$ws.Range("J34").Value = 102545.45
$ws.Range("K34").Value = 439.66668
$ws.Range("L34").Value = 307636.35
$ws.Range("M34").Value = -355.66668
$ws.Range("N34").Value = -307804.35
$ws.Range("H109").Value = 250723.5
$ws.Range("I109").Value = 964.6667
$ws.Range("J109").Value = 1000000
$ws.Range("K109").Value = 2894.0001
$ws.Range("L109").Value = 3000000
$ws.Range("M109").Value = -1854.0001
$ws.Range("N109").Value = -3002080
$ws.Range("H111").Value = 6463
$ws.Range("I111").Value = 6463
$ws.Range("K111").Value = 19389
$ws.Range("M111").Value = -16322
$ws.Range("H121").Value = 716227.3
$ws.Range("I121").Value = 1375
$ws.Range("K121").Value = 4125
$ws.Range("M121").Value = -2815
$ws.Range("H140").Value = 975.625
$ws.Range("I140").Value = 975.625
$ws.Range("K140").Value = 2926.875
$ws.Range("M140").Value = 2253.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8232.23
$ws.Range("I113").Value = 7335.1665
$ws.Range("K113").Value = 7335.1665
$ws.Range("M113").Value = -5165.1665
$ws.Range("H122").Value = 3499.9333
$ws.Range("I122").Value = 3350.3333
$ws.Range("K122").Value = 10050.9999
$ws.Range("M122").Value = -7600.999899999999
$ws.Range("H132").Value = 236789.52
$ws.Range("I132").Value = 246209.95
$ws.Range("K132").Value = 738629.8500000001
$ws.Range("M132").Value = -736099.8500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3305
$ws.Range("I22").Value = 3074.2
$ws.Range("J22").Value = 3433.2222
$ws.Range("K22").Value = 3074.2
$ws.Range("L22").Value = 3433.2222
$ws.Range("M22").Value = -2779.2
$ws.Range("N22").Value = -4023.2222
$ws.Range("H27").Value = 3305
$ws.Range("I27").Value = 3074.2
$ws.Range("J27").Value = 3433.2222
$ws.Range("K27").Value = 3074.2
$ws.Range("L27").Value = 3433.2222
$ws.Range("M27").Value = -2967.2
$ws.Range("N27").Value = -3647.2222
$ws.Range("H40").Value = 5467.3
$ws.Range("I40").Value = 2943.25
$ws.Range("K40").Value = 2943.25
$ws.Range("M40").Value = -2807.25
$ws.Range("H68").Value = 5600
$ws.Range("I68").Value = 5600
$ws.Range("K68").Value = 5600
$ws.Range("M68").Value = -4851
$ws.Range("H71").Value = 5600
$ws.Range("I71").Value = 5600
$ws.Range("K71").Value = 28000
$ws.Range("M71").Value = -24256
$ws.Range("H122").Value = 1832319.4
$ws.Range("I122").Value = 1687251.4
$ws.Range("J122").Value = 2006401
$ws.Range("K122").Value = 5061754.199999999
$ws.Range("L122").Value = 6019203
$ws.Range("M122").Value = -5059304.199999999
$ws.Range("N122").Value = -6024103
$ws.Range("H136").Value = 5407.7393
$ws.Range("I136").Value = 2747.2856
$ws.Range("J136").Value = 6571.6875
$ws.Range("K136").Value = 8241.856800000001
$ws.Range("L136").Value = 19715.0625
$ws.Range("M136").Value = -5691.856800000001
$ws.Range("N136").Value = -24815.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 79676.766
$ws.Range("I96").Value = 113666.664
$ws.Range("J96").Value = 3199.5
$ws.Range("K96").Value = 113666.664
$ws.Range("L96").Value = 3199.5
$ws.Range("M96").Value = -112293.664
$ws.Range("N96").Value = -5945.5
$ws.Range("H113").Value = 2895.3333
$ws.Range("I113").Value = 1462.3334
$ws.Range("J113").Value = 4328.3335
$ws.Range("K113").Value = 4387.0002
$ws.Range("L113").Value = 12985.0005
$ws.Range("M113").Value = -2217.0002
$ws.Range("N113").Value = -17325.0005
$ws.Range("H122").Value = 31253908
$ws.Range("I122").Value = 52634108
$ws.Range("K122").Value = 157902324
$ws.Range("M122").Value = -157899874
